# db.xlsx - "hoteles" sheet cleanup: drop the redundant numeric "hoteles" id
# column, rename the hotel-name column header from "titulo" to "nombre", and
# shift the address column left to follow it.
#
# Before: A=id | B=hoteles(1,1) | C=titulo(hotel santiago/vina) | D=direccion
# After:  A=id | B=nombre(hotel santiago/vina) | C=direccion

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hoteles")

# Remove the now-unused column B ("hoteles" header, numeric 1/1 values).
# This shifts old C (titulo/name) -> B and old D (direccion) -> C, and the
# widths/bestFit column formatting carried on those columns shifts with them.
$ws.Columns.Item(2).Delete()

# Relabel the (now) B column header from "titulo" to "nombre".
$ws.Range("B1").Value = "nombre"

# The header cells no longer carry the old bold-ish style index, matching
# the cleaned-up styling in the edited workbook.
$ws.Range("A1:B1").ClearFormats()

# Match the saved selection/cursor position from the edit.
$ws.Range("B2").Select() | Out-Null
